$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell references to their new text values (as they appear in the
# source "before" data export - plain text cells, not numeric).
$updates = @{
    "D2" = "307.01"
    "D3" = "39.80"
    "E3" = "-7.58%"
    "D4" = "5.123"
    "E4" = "-1.67%"
    "D5" = "0.07699"
    "E5" = "-6.01%"
    "D6" = "4.241"
    "E6" = "-1.86%"
    "D7" = "1.625"
    "E7" = "-11.98%"
    "D8" = "0.8794"
    "E8" = "-5.96%"
    "D9" = "0.09979"
    "E9" = "-10.47%"
    "E10" = "-6.45%"
    "D11" = "0.08920"
    "E11" = "-5.57%"
    "D12" = "0.04407"
    "E12" = "-4.68%"
    "E13" = "-0.28%"
    "D14" = "0.001258"
    "E14" = "-2.39%"
    "D15" = "0.005880"
    "E15" = "3.75%"
    "E18" = "-0.51%"
    "D19" = "7.018"
    "E19" = "-5.41%"
    "D20" = "0.1339"
    "E20" = "-3.49%"
    "E21" = "14.56%"
    "D22" = "0.04151"
    "E22" = "-0.30%"
    "D23" = "0.001199"
    "E23" = "-4.03%"
    "D24" = "0.004051"
    "E24" = "-5.94%"
    "D25" = "0.0001220"
    "E25" = "10.98%"
    "E26" = "0.16%"
    "D38" = "0.02331"
    "E38" = "-14.02%"
    "D39" = "0.05148"
    "E39" = "-7.28%"
    "D40" = "0.007924"
    "E40" = "-0.56%"
    "D41" = "0.1323"
    "E41" = "-5.11%"
    "D42" = "0.006324"
    "E42" = "-3.30%"
    "D43" = "0.001957"
    "E43" = "-6.43%"
    "D44" = "0.008489"
    "E44" = "13.84%"
    "D45" = "0.3062"
    "E45" = "-4.57%"
    "D46" = "0.00006511"
    "E46" = "-6.45%"
    "E48" = "98.62%"
    "D49" = "0.007135"
    "E49" = "106.07%"
    "D50" = "0.00002100"
    "D51" = "0.0002000"
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    # Leading apostrophe forces Excel to store the value as literal text,
    # preserving formatting such as trailing zeros (e.g. 0.0002000) instead
    # of collapsing it to a number.
    $cell.Value = "'" + $updates[$ref]
    # Reset style to Normal so no stray quote-prefix / text-format style
    # attribute is left attached to the cell (matches original formatting).
    $cell.Style = "Normal"
}
